$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("SMR20190813")
$src.Copy($src)
$ws = $wb.Worksheets.Item(1)
$ws.Name = "SMR20200323"

# Insert 3 rows before row 19 (room for C91,C92,C93)
$ws.Rows("19:21").Insert()

# Delete old R35 (row26) and R36 (row27) rows (no longer used)
$ws.Rows("26:27").Delete()

# Insert 4 rows before row 29 (room for R91,R92,R93,R94)
$ws.Rows("29:32").Insert()
